# Fixing geopoint in shared_table model
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# Insert a new column before column C: this shifts the old
# schema.elementType / latitude / longitude / altitude / accuracy
# columns one place to the right (C->D, D->E, E->F, F->G, G->H)
# and preserves their custom widths, matching the diff's col-width shift.
$ws.Columns("C").Insert()

# New header: schema.name (elementType repeats for the geopoint object)
$ws.Range("C1").Value = "schema.name"
$ws.Range("C4").Value = "geopoint"

# Rename the shifted geopoint-property headers to the ".type" suffixed
# versions.
$ws.Range("E1").Value = "schema.properties.latitude.type"
$ws.Range("F1").Value = "schema.properties.longitude.type"
$ws.Range("G1").Value = "schema.properties.altitude.type"
$ws.Range("H1").Value = "schema.properties.accuracy.type"
